$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.109.48"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "2.493.61"
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'318.51"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("D6").Value = "'105.25"
$ws.Range("E6").Value = "  -3.68%  "

$ws.Range("D7").Value = "'0.520"
$ws.Range("E7").Value = "  -1.84%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'0.537"
$ws.Range("E9").Value = "  -1.77%  "

$ws.Range("D10").Value = "'38.84"
$ws.Range("E10").Value = "  -2.87%  "

$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("D12").Value = "'0.0800"
$ws.Range("E12").Value = "  -2.18%  "

$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("E14").Value = "  -1.96%  "

$ws.Range("D15").Value = "2.882.14"
$ws.Range("E15").Value = "  -0.76%  "

$ws.Range("D16").Value = "2.451.15"
$ws.Range("E16").Value = "  -2.34%  "

$ws.Range("D17").Value = "'0.838"
$ws.Range("E17").Value = "  -1.13%  "

$ws.Range("D18").Value = "47.994.06"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").Value = "'12.78"
$ws.Range("E19").Value = "  -3.24%  "

$ws.Range("D20").Value = "'2.92"
$ws.Range("E20").Value = "  +7.86%  "

$ws.Range("D21").Value = "'6.55"
$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").Value = "0.0₃0930"
$ws.Range("E22").Value = "  -1.53%  "

$ws.Range("D23").Value = "'280.78"
$ws.Range("E23").Value = "  +1.97%  "

$ws.Range("D24").Value = "'70.84"
$ws.Range("E24").Value = "  -1.49%  "

$ws.Range("D25").Value = "'2.50"
$ws.Range("E25").Value = "  -2.61%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").Value = "'25.68"
$ws.Range("E27").Value = "  -0.86%  "

$ws.Range("E28").Value = "  -7.90%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.139"
$ws.Range("E29").Value = "  -1.40%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'9.58"
$ws.Range("E30").Value = "  -4.92%  "

$ws.Range("D31").Value = "'34.65"
$ws.Range("E31").Value = "  -2.34%  "

$ws.Range("D32").Value = "'49.00"
$ws.Range("E32").Value = "  -0.96%  "

$ws.Range("D33").Value = "'19.29"
$ws.Range("E33").Value = "  -1.03%  "

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("E35").Value = "  -1.79%  "

$ws.Range("D36").Value = "'0.0771"
$ws.Range("E36").Value = "  -1.43%  "

$ws.Range("E37").Value = "  -0.88%  "

$ws.Range("D38").Value = "'4.50"
$ws.Range("E38").Value = "  -2.82%  "

$ws.Range("D39").Value = "'2.88"
$ws.Range("E39").Value = "  -2.61%  "

$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "'2.20"
$ws.Range("E41").Value = "  -0.43%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'118.92"
$ws.Range("E42").Value = "  -2.31%  "

$ws.Range("D43").Value = "'21.46"
$ws.Range("E43").Value = "  -1.73%  "

$ws.Range("D44").Value = "'0.0299"
$ws.Range("E44").Value = "  -2.37%  "

$ws.Range("D45").Value = "1.988.82"
$ws.Range("E45").Value = "  -1.89%  "

$ws.Range("D46").Value = "'3.13"
$ws.Range("E46").Value = "  +0.66%  "

$ws.Range("E47").Value = "  +6.02%  "

$ws.Range("D48").Value = "'1.94"
$ws.Range("E48").Value = "  +4.29%  "

$ws.Range("D49").Value = "'8.96"
$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("D50").Value = "'5.11"
$ws.Range("E50").Value = "  -1.26%  "

$ws.Range("D51").Value = "'79.33"
